$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value2 = 356
$ws.Cells.Item(7, 6).Value2 = 1166
$ws.Cells.Item(9, 6).Value2 = 7097
$ws.Cells.Item(12, 6).Value2 = 2042
$ws.Cells.Item(13, 6).Value2 = 7968
$ws.Cells.Item(16, 6).Value2 = 5504
$ws.Cells.Item(19, 6).Value2 = 1022
$ws.Cells.Item(20, 6).Value2 = 4562
$ws.Cells.Item(22, 6).Value2 = 384
$ws.Cells.Item(25, 6).Value2 = 374
$ws.Cells.Item(28, 6).Value2 = 2327
$ws.Cells.Item(31, 6).Value2 = 77
$ws.Cells.Item(32, 6).Value2 = 139
$ws.Cells.Item(33, 6).Value2 = 577
$ws.Cells.Item(36, 6).Value2 = 1485
$ws.Cells.Item(39, 6).Value2 = 2303
$ws.Cells.Item(40, 6).Value2 = 2211

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value2 = 89
$ws.Cells.Item(4, 6).Value2 = 60
$ws.Cells.Item(6, 6).Value2 = 24

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value2 = 1278

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value2 = 1278
$ws.Cells.Item(6, 6).Value2 = 89
$ws.Cells.Item(7, 6).Value2 = 356
$ws.Cells.Item(8, 6).Value2 = 1166
$ws.Cells.Item(10, 6).Value2 = 7097
$ws.Cells.Item(13, 6).Value2 = 2042
$ws.Cells.Item(14, 6).Value2 = 7968
$ws.Cells.Item(17, 6).Value2 = 5504
$ws.Cells.Item(20, 6).Value2 = 1022
$ws.Cells.Item(21, 6).Value2 = 4562
$ws.Cells.Item(23, 6).Value2 = 384
$ws.Cells.Item(27, 6).Value2 = 60
$ws.Cells.Item(28, 6).Value2 = 374
$ws.Cells.Item(30, 6).Value2 = 2327
$ws.Cells.Item(33, 6).Value2 = 77
$ws.Cells.Item(34, 6).Value2 = 139
$ws.Cells.Item(36, 6).Value2 = 577
$ws.Cells.Item(39, 6).Value2 = 24
$ws.Cells.Item(40, 6).Value2 = 1485
$ws.Cells.Item(43, 6).Value2 = 2303
$ws.Cells.Item(45, 6).Value2 = 2211
